$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) format, used to restore style
# after forcing a cell to Text format so numeric-looking strings are not coerced.
$blank = $ws.Range("D9")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.015.41"
$ws.Range("D2").Style = $blank.Style
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.088.96"
$ws.Range("D3").Style = $blank.Style
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.02"
$ws.Range("D5").Style = $blank.Style
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.95"
$ws.Range("D6").Style = $blank.Style
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.084.93"
$ws.Range("D8").Style = $blank.Style
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.45"
$ws.Range("D11").Style = $blank.Style
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("E13").Value = "  +4.96%  "
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.595.79"
$ws.Range("D15").Style = $blank.Style
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.042.19"
$ws.Range("D16").Style = $blank.Style
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.090.61"
$ws.Range("D18").Style = $blank.Style
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.69"
$ws.Range("D19").Style = $blank.Style
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "484.93"
$ws.Range("D20").Style = $blank.Style
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.40"
$ws.Range("D21").Style = $blank.Style
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.701"
$ws.Range("D22").Style = $blank.Style
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.12"
$ws.Range("D23").Style = $blank.Style
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.77"
$ws.Range("D24").Style = $blank.Style
$ws.Range("E24").Value = "  +2.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.23"
$ws.Range("D25").Style = $blank.Style
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.10"
$ws.Range("D28").Style = $blank.Style
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.43"
$ws.Range("D30").Style = $blank.Style
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.90"
$ws.Range("D31").Style = $blank.Style
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.15"
$ws.Range("D32").Style = $blank.Style
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "57.65"
$ws.Range("D33").Style = $blank.Style
$ws.Range("E33").Value = "  -3.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.35"
$ws.Range("D34").Style = $blank.Style
$ws.Range("E34").Value = "  -6.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "502.57"
$ws.Range("D35").Style = $blank.Style
$ws.Range("E35").Value = "  -5.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.36"
$ws.Range("D36").Style = $blank.Style
$ws.Range("E36").Value = "  +3.64%  "
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.239.28"
$ws.Range("D38").Style = $blank.Style
$ws.Range("E38").Value = "  +5.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0398"
$ws.Range("D39").Style = $blank.Style
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0796"
$ws.Range("D40").Style = $blank.Style
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.13"
$ws.Range("D42").Style = $blank.Style
$ws.Range("E42").Value = "  +0.42%  "
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.255"
$ws.Range("D44").Style = $blank.Style
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "122.80"
$ws.Range("D46").Style = $blank.Style
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.05"
$ws.Range("D47").Style = $blank.Style
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₃0532"
$ws.Range("D48").Style = $blank.Style
$ws.Range("E48").Value = "  +5.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.66"
$ws.Range("D49").Style = $blank.Style
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("E51").Value = "  +2.12%  "
